$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The rule table (rows 18-26) used to have 3 ACTION columns (C, D, E).
# Collapse it down to a single ACTION column: drop the old column C
# (its values are being replaced by what used to be in column D), then
# shift column D left into C, and drop the (now stray) trailing column
# that used to be E. Net effect: old column D becomes the new column C,
# and columns D & E disappear entirely.
$ws.Range("C18:C26").Delete(-4159)   # xlShiftToLeft: drop old C, D->C, E->D
$ws.Range("D18:D26").ClearContents() # drop the now-orphaned former column E

# Row 26 didn't have any data in the old D/E columns, but the new
# single ACTION column for that row should read "Test".
$ws.Cells.Item(26, 3).Value2 = "Test"

$ws.Range("A1:C26").Select() | Out-Null
